$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("H10").ClearFormats()
$ws.Range("H10").HorizontalAlignment = -4108
$ws.Range("H10").VerticalAlignment = -4108
$ws.Range("H10").Orientation = 45
$ws.Range("H10").Font.Name = "Aptos Narrow"
$ws.Range("H10").Font.Bold = $true
$ws.Range("H10").Font.Size = 16
